# Apply changes described by the diff:
#  - Fill column D (survival rate) values for rows 15-90 in Sheet1
#  - Update the active selection / scroll position of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ranges of rows (inclusive) mapped to the value that should be written
# into column D for every row in that range.
$ranges = @(
    @{Start = 15; End = 39; Value = 84.9},
    @{Start = 40; End = 49; Value = 90},
    @{Start = 50; End = 59; Value = 91.2},
    @{Start = 60; End = 69; Value = 92.4},
    @{Start = 70; End = 79; Value = 83},
    @{Start = 80; End = 90; Value = 70.3}
)

foreach ($r in $ranges) {
    $ws.Range("D$($r.Start):D$($r.End)").Value = $r.Value
}

# Scroll the sheet so row 62 is the top visible row and select F89,
# matching the saved view state captured in the workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 62
$ws.Range("F89").Select()
